$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 888
$ws1.Range("F3").Value = 1012
$ws1.Range("F7").Value = 689
$ws1.Range("F8").Value = 160
$ws1.Range("F9").Value = 1291
$ws1.Range("F12").Value = 547
$ws1.Range("F13").Value = 183
$ws1.Range("F14").Value = 40
$ws1.Range("F15").Value = 1001
$ws1.Range("F16").Value = 19
$ws1.Range("F17").Value = 410
$ws1.Range("F18").Value = 375
$ws1.Range("F19").Value = 93
$ws1.Range("F20").Value = 590
$ws1.Range("F21").Value = 150
$ws1.Range("F22").Value = 639
$ws1.Range("F24").Value = 1018
$ws1.Range("F25").Value = 16

# Sheet "演出" (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 245
$ws2.Range("F11").Value = 113

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 888
$ws4.Range("F5").Value = 1012
$ws4.Range("F9").Value = 689
$ws4.Range("F10").Value = 160
$ws4.Range("F11").Value = 1291
$ws4.Range("F16").Value = 547
$ws4.Range("F18").Value = 183
$ws4.Range("F19").Value = 40
$ws4.Range("F20").Value = 1001
$ws4.Range("F22").Value = 19
$ws4.Range("F23").Value = 410
$ws4.Range("F24").Value = 375
$ws4.Range("F25").Value = 93
$ws4.Range("F26").Value = 245
$ws4.Range("F28").Value = 590
$ws4.Range("F31").Value = 113
$ws4.Range("F32").Value = 113
$ws4.Range("F33").Value = 150
$ws4.Range("F34").Value = 639
$ws4.Range("F36").Value = 1018
$ws4.Range("F37").Value = 16

$wb.Save()
